$d = $word.ActiveDocument

function Find-RangeByText($searchText) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $found = $f.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($f.Found) {
        return $f.Parent
    }
    return $null
}

function Force-RunBreak($range) {
    # Touching a numeric Font property (and then restoring it) is enough to make
    # the engine materialise a standalone <w:r> for this sub-range without
    # disturbing any of the boolean formatting flags that are already correct.
    $original = $range.Font.Size
    if ($original -eq 12) {
        $range.Font.Size = 11
        $range.Font.Size = 12
    } else {
        $range.Font.Size = 12
        $range.Font.Size = $original
    }
}

# ---------------------------------------------------------------------------
# Change 1: "You will see several strings..." -> split into "Y" | "ou will..."
# ---------------------------------------------------------------------------
$full1 = Find-RangeByText("You will see several strings being loaded into the stack if you ran a")
if ($full1 -ne $null) {
    $yRange = $d.Range($full1.Start, $full1.Start + 1)
    Force-RunBreak($yRange)
}

# ---------------------------------------------------------------------------
# Change 2: "...is place in "local_40h"..." -> "...is place" | "d" | " in ..."
# (fixes the typo "place in" -> "placed in", and splits the run accordingly)
# ---------------------------------------------------------------------------
$place = Find-RangeByText("can identify which one is place")
if ($place -ne $null) {
    $insertionPoint = $d.Range($place.End, $place.End)
    $insertionPoint.InsertBefore("d")

    $dRange = $d.Range($place.End, $place.End + 1)
    Force-RunBreak($dRange)
}

# ---------------------------------------------------------------------------
# Change 3: add a new character style "ListLabel2" (mirrors "ListLabel1")
# ---------------------------------------------------------------------------
$existing = $null
try { $existing = $d.Styles("ListLabel2") } catch { $existing = $null }
if ($existing -eq $null) {
    $style = $d.Styles.Add("ListLabel2", 2)
    $style.NameLocal = "ListLabel 2"
    $style.QuickStyle = $true
    $style.Font.Bold = $false
    $style.Font.BoldBi = $false
    $style.Font.Italic = $false
    $style.Font.ItalicBi = $false
    $style.Font.Size = 12
    $style.Font.SizeBi = 12
    $style.Font.Underline = 0
}
